$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.394.95"
$ws.Range("E2").Value = "  +0.07%  "
$ws.Range("D3").Value = "1.848.16"
$ws.Range("E3").Value = "  +0.02%  "
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'240.56"
$ws.Range("E5").Value = "  +0.20%  "
$ws.Range("D6").Value = "'0.6266"
$ws.Range("E6").Value = "  -0.30%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "'0.07493"
$ws.Range("E8").Value = "  -1.79%  "
$ws.Range("D9").Value = "'0.2904"
$ws.Range("E9").Value = "  -0.08%  "
$ws.Range("D10").Value = "'24.39"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.848.06"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "'5.001"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("D14").Value = "'0.6804"
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("D15").Value = "'0.00001046"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").Value = "'82.20"
$ws.Range("E16").Value = "  -1.24%  "
$ws.Range("D17").Value = "2.103.48"
$ws.Range("E17").Value = "  -3.82%  "
$ws.Range("D18").Value = "'6.174"
$ws.Range("E18").Value = "  +0.38%  "
$ws.Range("D19").Value = "29.429.62"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("D20").Value = "'229.67"
$ws.Range("E20").Value = "  +0.93%  "
$ws.Range("D21").Value = "'12.34"
$ws.Range("E21").Value = "  -0.03%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "'7.458"
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").Value = "'158.78"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").Value = "'0.1375"
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("D27").Value = "'8.416"
$ws.Range("E27").Value = "  +0.22%  "
$ws.Range("D28").Value = "'17.56"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("D29").Value = "'0.06437"
$ws.Range("E29").Value = "  +14.99%  "
$ws.Range("D30").Value = "'1.392"
$ws.Range("E30").Value = "  +1.42%  "
$ws.Range("D31").Value = "'1.478"
$ws.Range("E31").Value = "  +1.35%  "
$ws.Range("D32").Value = "'4.097"
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("D33").Value = "'4.092"
$ws.Range("E33").Value = "  +0.41%  "
$ws.Range("D34").Value = "'1.829"
$ws.Range("E34").Value = "  -0.37%  "
$ws.Range("D35").Value = "'1.142"
$ws.Range("E35").Value = "  -1.70%  "
$ws.Range("D36").Value = "'0.6995"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("D37").Value = "'2.580"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "1.266.00"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").Value = "'2.833"
$ws.Range("E39").Value = "  +4.53%  "
$ws.Range("D40").Value = "'0.01828"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("D41").Value = "'6.607"
$ws.Range("E41").Value = "  +3.43%  "
$ws.Range("D42").Value = "'0.9100"
$ws.Range("E42").Value = "  +0.70%  "
$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "2.008.71"
$ws.Range("E44").Value = "  -18.37%  "
$ws.Range("D45").Value = "'101.55"
$ws.Range("E45").Value = "  +0.07%  "
$ws.Range("D46").Value = "'66.35"
$ws.Range("E46").Value = "  +0.42%  "
$ws.Range("D47").Value = "'1.752"
$ws.Range("E47").Value = "  +4.82%  "
$ws.Range("D48").Value = "'7.080"
$ws.Range("E48").Value = "  -1.35%  "
$ws.Range("D49").Value = "'0.1174"
$ws.Range("E49").Value = "  +3.21%  "
$ws.Range("D50").Value = "'9.021"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'0.3950"
$ws.Range("E51").Value = "  -1.45%  "
